$d = $word.ActiveDocument

# --- Change 1: "On how it calculates load on the machine beams:"
#     -> "On how it calculates" | " loads" | ":"  (three runs, same rPr)
$p1 = $d.Paragraphs(1)
$full = $p1.Range
$start = $full.Start
$text = $full.Text
$prefix = "On how it calculates"
$idx = $prefix.Length
$colonIdx = $text.IndexOf(":")
$midRange = $d.Range($start + $idx, $start + $colonIdx)
$midRange.Text = " loads"
# Toggling a character property on just the replaced span and back forces
# Word to keep it as a distinct run instead of re-merging it with its
# now-identically-formatted neighbours.
$newRange = $d.Range($start + $idx, $start + $idx + 6)
$newRange.Bold = 1
$newRange.Bold = 0

# --- Change 2: remove the stray empty paragraph between
#     "...at the machine beams." and "The point is located..."
# (A paragraph's Range.Text is just the paragraph mark "\r" when the
#  paragraph itself has no runs, so compare against that rather than "".)
$cr = [string][char]13
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq $cr -and $p.Range.InlineShapes.Count -eq 0) {
        $prevText = ""
        if ($i -gt 1) { $prevText = $d.Paragraphs($i - 1).Range.Text }
        if ($prevText.TrimEnd([char]13).EndsWith("at the machine beams.")) {
            $p.Range.Delete()
            break
        }
    }
}

# --- Change 3: mark the run holding the illustration <w:drawing> as NoProof
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.NoProofing = 1
    }
}
